# Commit: "add cellular data models"
# The deck gains a second, blank slide (placed right after the existing
# slide) that will later hold the new cellular-data-model content.

$p = $ppt.ActivePresentation

# 12 == ppLayoutBlank: insert a brand-new blank slide at position 2.
$newSlide = $p.Slides.Add(2, 12)
